$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 920
$ws1.Range("F3").Value = 1484
$ws1.Range("F4").Value = 1156
$ws1.Range("F5").Value = 542
$ws1.Range("F7").Value = 10
$ws1.Range("F8").Value = 699
$ws1.Range("F9").Value = 305
$ws1.Range("F11").Value = 108
$ws1.Range("F12").Value = 229
$ws1.Range("F13").Value = 177
$ws1.Range("F14").Value = 177
$ws1.Range("F15").Value = 4093
$ws1.Range("F16").Value = 26
$ws1.Range("F18").Value = 451
$ws1.Range("F20").Value = 518
$ws1.Range("F21").Value = 300
$ws1.Range("F24").Value = 10
$ws1.Range("F25").Value = 686
$ws1.Range("F30").Value = 1645
$ws1.Range("F31").Value = 375

# Sheet: 演出 (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 247

# Sheet: 本地生活 (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 131

# Sheet: 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 920
$ws4.Range("F4").Value = 1484
$ws4.Range("F5").Value = 1156
$ws4.Range("F8").Value = 131
$ws4.Range("F9").Value = 542
$ws4.Range("F11").Value = 10
$ws4.Range("F12").Value = 699
$ws4.Range("F14").Value = 305
$ws4.Range("F16").Value = 108
$ws4.Range("F17").Value = 229
$ws4.Range("F18").Value = 177
$ws4.Range("F19").Value = 177
$ws4.Range("F20").Value = 4094
$ws4.Range("F21").Value = 26
$ws4.Range("F24").Value = 451
$ws4.Range("F26").Value = 518
$ws4.Range("F27").Value = 300
$ws4.Range("F31").Value = 10
$ws4.Range("F32").Value = 247
$ws4.Range("F35").Value = 686
$ws4.Range("F43").Value = 1645
$ws4.Range("F44").Value = 375
